$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "id" column before the current column B (nombre) ---
$ws.Columns.Item(2).Insert()

# The insert leaves the new B1 without the header style, fix it up by
# copying the (shifted) neighbour's format, then set its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "id"

# The insert also left B2:B3 carrying the bordered/bold header style
# inherited from the insert op - the data column should be unstyled.
$ws.Range("B2:B3").ClearFormats()

# --- Update row 2 (existing user, now gets an id + new data) ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Juan"
$ws.Range("D2").Value = "Penas"
$ws.Range("E2").Value = "TipoUsuario.Ventas"
$ws.Range("F2").Value = 1

# --- Update row 3 (existing user, now gets an id + new data) ---
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Alberto"
$ws.Range("D3").Value = "Erreguin"
$ws.Range("E3").Value = "TipoUsuario.Inventario"
$ws.Range("F3").Value = 1

# --- New row 4 ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "John"
$ws.Range("D4").Value = "Doe"
$ws.Range("F4").Value = 1

# --- New row 5 ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Martha"
$ws.Range("D5").Value = "Almaraz"
$ws.Range("E5").Value = "TipoUsuario.Administrador"
$ws.Range("F5").Value = 1
